# Update "想去人数" (F) and "最低票价" (G) figures to the latest scraped
# values for the gh-pages data refresh (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Changes that apply identically to both the "展览" and "全部类型" sheets.
$commonChanges = @{
    "F3"  = 1471
    "F4"  = 183
    "F6"  = 258
    "G6"  = 0
    "F7"  = 110
    "F9"  = 204
    "F12" = 4925
    "F14" = 7204
    "F20" = 15
    "F22" = 1694
    "F23" = 93
    "F24" = 88
    "F25" = 2816
    "F31" = 426
    "F32" = 265
    "F33" = 71
    "F35" = 1120
    "F37" = 1320
    "F38" = 98
    "F43" = 40
    "F45" = 2572
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($ref in $commonChanges.Keys) {
        $ws.Range($ref).Value = $commonChanges[$ref]
    }
}

# "全部类型" has one additional change not present on "展览".
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F39").Value = 0
